$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.333330166666667
$ws.Range("A3").Value = 0.333330166666667
$ws.Range("A4").Value = 0.333330166666667
$ws.Range("A5").Value = 0.333330166666667
$ws.Range("A6").Value = 1.49998575
